# Generate Report for Handback
# Updates the "Xliff Generate Date" / handoff-handback timestamp cells to
# reflect a fresh report generation run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-06 11:06:51"

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-06 11:06:40"
$wsZhCn.Range("K2").Value = "2016-09-06 11:07:37"

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-06 11:06:51"
$wsDeDe.Range("K2").Value = "2016-09-06 11:07:57"
